## Minor correction on the "Maintaining Context During Parsing" slide:
## Split the paragraph describing "exit statement" semantics into three
## runs, fixing the stray period in "a loop.," -> "a loop,".
##
## Before (single run):
##   "An exit statement has meaning only when nested inside a loop., and
##    code generation for an exit statement requires knowledge of which
##    loop encloses it."
##
## After (three runs):
##   [1] "An exit statement has meaning only when nested inside "
##   [2] "a loop, "
##   [3] "and code generation for an exit statement requires knowledge of
##        which loop encloses it."

$p = $ppt.ActivePresentation

$oldSnippet = "An exit statement has meaning only when nested inside a loop"

$run1 = "An exit statement has meaning only when nested inside "
$run2 = "a loop, "
$run3 = "and code generation for an exit statement requires knowledge of which loop encloses it."
$newText = $run1 + $run2 + $run3

# Locate the paragraph holding the text, searching every slide/shape so the
# script is resilient to any reordering.
$targetSlide = $null
$targetShape = $null
$targetParaIndex = 0

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $textRange = $shape.TextFrame.TextRange
            $paraCount = $textRange.Paragraphs().Count
            for ($pi = 1; $pi -le $paraCount; $pi++) {
                $para = $textRange.Paragraphs($pi)
                if ($para.Text.StartsWith($oldSnippet)) {
                    $targetSlide = $slide
                    $targetShape = $shape
                    $targetParaIndex = $pi
                }
            }
        }
    }
}

if ($targetShape -ne $null) {
    $textRange = $targetShape.TextFrame.TextRange
    $para = $textRange.Paragraphs($targetParaIndex)
    $paraStart = $para.Start

    # Replace the whole paragraph's text with the corrected wording (single
    # run for now), then carve out the middle run "a loop, " so it ends up
    # as its own <a:r> -- matching the three-run structure in the target.
    $para.Text = $newText

    $midStart = $paraStart + $run1.Length
    $midRange = $textRange.Characters($midStart, $run2.Length)
    $midRange.Text = $midRange.Text
}
